$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update timestamp and swap country-name labels (shared-string reorder effect)
$ws.Range("A1").Value = "Datos actualizados a 10 de Septiembre de 2020 a las 08:26"
$ws.Range("A62").Value = "Uzbekistan"
$ws.Range("A63").Value = "Armenia"
$ws.Range("A152").Value = "Georgia"
$ws.Range("A153").Value = "Nueva Zelanda"
$ws.Range("A214").Value = "Montserrat"
$ws.Range("A215").Value = "Islas Malvinas"

# Update statistic values (refreshed COVID data + re-ranked rows)
$ws.Range("B27").Value = 145612
$ws.Range("C27").Value = 2582
$ws.Range("D27").Value = 65877
$ws.Range("E27").Value = 76712
$ws.Range("G27").Value = 44
$ws.Range("H27").Value = 3023
$ws.Range("B28").Value = 142582
$ws.Range("C28").Value = 1485
$ws.Range("D28").Value = 109757
$ws.Range("E28").Value = 31771
$ws.Range("B62").Value = 45160
$ws.Range("C62").Value = 230
$ws.Range("D62").Value = 42212
$ws.Range("E62").Value = 2580
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = 368
$ws.Range("B63").Value = 45152
$ws.Range("D63").Value = 41023
$ws.Range("E63").Value = 3224
$ws.Range("H63").Value = 905
$ws.Range("D74").Value = 16807
$ws.Range("E74").Value = 9021
$ws.Range("G74").Value = 4
$ws.Range("H74").Value = 774
$ws.Range("B75").Value = 26513
$ws.Range("C75").Value = 48
$ws.Range("D75").Value = 23036
$ws.Range("E75").Value = 2689
$ws.Range("B152").Value = 1830
$ws.Range("C152").Value = 57
$ws.Range("D152").Value = 1334
$ws.Range("E152").Value = 477
$ws.Range("H152").Value = 19
$ws.Range("B153").Value = 1792
$ws.Range("C153").Value = 4
$ws.Range("D153").Value = 1648
$ws.Range("E153").Value = 120
$ws.Range("H153").Value = 24
$ws.Range("B176").Value = 496
$ws.Range("C176").Value = 1
$ws.Range("E176").Value = 14
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0
